$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 94.5
$ws.Range("I12").Value = 94.5
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 94.5
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 75.5
$ws.Range("N12").ClearContents()
$ws.Range("H74").Value = 4183
$ws.Range("I74").Value = 4183
$ws.Range("K74").Value = 4183
$ws.Range("M74").Value = -3247
$ws.Range("H77").Value = 4183
$ws.Range("I77").Value = 4183
$ws.Range("K77").Value = 20915
$ws.Range("M77").Value = -16235
$ws.Range("H88").Value = 1499
$ws.Range("I88").Value = 1333
$ws.Range("K88").Value = 1333
$ws.Range("M88").Value = -927
$ws.Range("H91").Value = 1499
$ws.Range("I91").Value = 1333
$ws.Range("K91").Value = 1333
$ws.Range("M91").Value = 71
$ws.Range("H107").Value = 722.4286
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H137").Value = 3677.4285
$ws.Range("J137").Value = 4090
$ws.Range("L137").Value = 12270
$ws.Range("N137").Value = -17370
$ws.Range("H141").Value = 4266.4546
$ws.Range("I141").Value = 4493.2
$ws.Range("K141").Value = 13479.6
$ws.Range("M141").Value = -8299.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 698.5
$ws.Range("I4").Value = 698.5
$ws.Range("K4").Value = 698.5
$ws.Range("M4").Value = -582.5
$ws.Range("H5").Value = 330.33334
$ws.Range("I5").Value = 330.33334
$ws.Range("K5").Value = 330.33334
$ws.Range("M5").Value = -218.33334
$ws.Range("H28").Value = 15568.25
$ws.Range("I28").Value = 4326
$ws.Range("K28").Value = 4326
$ws.Range("M28").Value = -4134
$ws.Range("H45").Value = 1750.1875
$ws.Range("I45").Value = 1580.3
$ws.Range("K45").Value = 1580.3
$ws.Range("M45").Value = -1203.3
$ws.Range("H99").Value = 15568.25
$ws.Range("I99").Value = 4326
$ws.Range("K99").Value = 4326
$ws.Range("M99").Value = -1331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 330.33334
$ws.Range("I4").Value = 330.33334
$ws.Range("K4").Value = 330.33334
$ws.Range("M4").Value = -215.33334
$ws.Range("H86").Value = 2798
$ws.Range("I86").Value = 2330
$ws.Range("K86").Value = 2330
$ws.Range("M86").Value = -1207
$ws.Range("H89").Value = 2798
$ws.Range("I89").Value = 2330
$ws.Range("K89").Value = 11650
$ws.Range("M89").Value = -6034
$ws.Range("H94").Value = 3731.6667
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 5000
$ws.Range("N94").Value = -5902
$ws.Range("H105").Value = 2485
$ws.Range("I105").Value = 1982
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 1982
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -235
$ws.Range("N105").Value = -8494
$ws.Range("H137").Value = 34999.5
$ws.Range("I137").Value = 34999.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 34999.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -29899.5
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 57.8
$ws.Range("I7").Value = 61.25
$ws.Range("J7").Value = 44
$ws.Range("K7").Value = 61.25
$ws.Range("L7").Value = 44
$ws.Range("M7").Value = 51.75
$ws.Range("N7").Value = -270
$ws.Range("H22").Value = 117
$ws.Range("I22").Value = 117
$ws.Range("K22").Value = 117
$ws.Range("M22").Value = 233
$ws.Range("H92").Value = 46820.2
$ws.Range("J92").Value = 51150.25
$ws.Range("L92").Value = 51150.25
$ws.Range("N92").Value = -56142.25
$ws.Range("H93").Value = 17266.334
$ws.Range("J93").Value = 16800
$ws.Range("L93").Value = 16800
$ws.Range("N93").Value = -20544
$ws.Range("H122").Value = 655.8333
$ws.Range("I122").Value = 629
$ws.Range("K122").Value = 1887
$ws.Range("M122").Value = 563

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 339.33334
$ws.Range("I2").Value = 8.25
$ws.Range("J2").Value = 1001.5
$ws.Range("K2").Value = 49.5
$ws.Range("L2").Value = 6009
$ws.Range("M2").Value = 63.5
$ws.Range("N2").Value = -6235
$ws.Range("H5").Value = 3824.5
$ws.Range("J5").Value = 2249
$ws.Range("L5").Value = 6747
$ws.Range("N5").Value = -6971
$ws.Range("H6").Value = 8783.571
$ws.Range("J6").Value = 20116.334
$ws.Range("L6").Value = 60349.00199999999
$ws.Range("N6").Value = -60575.00199999999
$ws.Range("H23").Value = 441.35715
$ws.Range("I23").Value = 242
$ws.Range("J23").Value = 707.1667
$ws.Range("K23").Value = 726
$ws.Range("L23").Value = 2121.5001
$ws.Range("M23").Value = -491
$ws.Range("N23").Value = -2591.5001
$ws.Range("H37").Value = 99997.5
$ws.Range("J37").Value = 99997.5
$ws.Range("L37").Value = 299992.5
$ws.Range("N37").Value = -300216.5
$ws.Range("H38").Value = 209.14285
$ws.Range("I38").Value = 203
$ws.Range("K38").Value = 609
$ws.Range("M38").Value = -262
$ws.Range("H75").Value = 3188
$ws.Range("I75").Value = 2774.5
$ws.Range("J75").Value = 4015
$ws.Range("K75").Value = 8323.5
$ws.Range("L75").Value = 12045
$ws.Range("M75").Value = -7325.5
$ws.Range("N75").Value = -14041
$ws.Range("H78").Value = 3188
$ws.Range("I78").Value = 2774.5
$ws.Range("J78").Value = 4015
$ws.Range("K78").Value = 24970.5
$ws.Range("L78").Value = 36135
$ws.Range("M78").Value = -19978.5
$ws.Range("N78").Value = -46119
$ws.Range("H92").Value = 502.5
$ws.Range("I92").Value = 502.5
$ws.Range("K92").Value = 1507.5
$ws.Range("M92").Value = -259.5
$ws.Range("H135").Value = 3824.5
$ws.Range("J135").Value = 2249
$ws.Range("L135").Value = 20241
$ws.Range("N135").Value = -25311

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 112.666664
$ws.Range("I2").Value = 80.625
$ws.Range("J2").Value = 176.75
$ws.Range("K2").Value = 80.625
$ws.Range("L2").Value = 176.75
$ws.Range("M2").Value = 32.375
$ws.Range("N2").Value = -402.75
$ws.Range("H6").Value = 950
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("H16").Value = 950
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H102").Value = 517.2
$ws.Range("I102").Value = 530.2222
$ws.Range("K102").Value = 530.2222
$ws.Range("M102").Value = 1091.7778
$ws.Range("H126").Value = 5253.909
$ws.Range("I126").Value = 5449.25
$ws.Range("J126").Value = 4733
$ws.Range("K126").Value = 16347.75
$ws.Range("L126").Value = 14199
$ws.Range("M126").Value = -13877.75
$ws.Range("N126").Value = -19139
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 327
$ws.Range("I16").Value = 327
$ws.Range("K16").Value = 327
$ws.Range("M16").Value = -157
$ws.Range("H22").Value = 1926.7273
$ws.Range("I22").Value = 1919.2
$ws.Range("J22").Value = 2002
$ws.Range("K22").Value = 1919.2
$ws.Range("L22").Value = 2002
$ws.Range("M22").Value = -1624.2
$ws.Range("N22").Value = -2592
$ws.Range("H27").Value = 1926.7273
$ws.Range("I27").Value = 1919.2
$ws.Range("J27").Value = 2002
$ws.Range("K27").Value = 1919.2
$ws.Range("L27").Value = 2002
$ws.Range("M27").Value = -1812.2
$ws.Range("N27").Value = -2216
$ws.Range("H40").Value = 2331
$ws.Range("I40").Value = 2331
$ws.Range("K40").Value = 2331
$ws.Range("M40").Value = -2195
$ws.Range("H55").Value = 302.6875
$ws.Range("J55").Value = 393.66666
$ws.Range("L55").Value = 393.66666
$ws.Range("N55").Value = -739.66666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 4951
$ws.Range("J4").Value = 5938.25
$ws.Range("L4").Value = 5938.25
$ws.Range("N4").Value = -6164.25
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H45").Value = 18500
$ws.Range("J45").Value = 18500
$ws.Range("L45").Value = 18500
$ws.Range("N45").Value = -19482
$ws.Range("H55").Value = 13682.167
$ws.Range("I55").Value = 600
$ws.Range("J55").Value = 16298.6
$ws.Range("K55").Value = 600
$ws.Range("L55").Value = 16298.6
$ws.Range("M55").Value = -323
$ws.Range("N55").Value = -16852.6
$ws.Range("H62").Value = 3380.7693
$ws.Range("I62").Value = 3038.6667
$ws.Range("K62").Value = 3038.6667
$ws.Range("M62").Value = -2414.6667
$ws.Range("H65").Value = 3380.7693
$ws.Range("I65").Value = 3038.6667
$ws.Range("K65").Value = 15193.3335
$ws.Range("M65").Value = -12073.3335
$ws.Range("H136").Value = 2363.9167
$ws.Range("I136").Value = 1808
$ws.Range("J136").Value = 4031.6667
$ws.Range("K136").Value = 5424
$ws.Range("L136").Value = 12095.0001
$ws.Range("M136").Value = -2874
$ws.Range("N136").Value = -17195.0001
